# Auto-generated edit script applying numeric updates to the Leve profit sheets
# per the commit diff (Marilith_Profits.xlsx -> Sheets/*).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 518.6774
$ws.Range("I15").Value = 518.6774
$ws.Range("K15").Value = 1556.0322
$ws.Range("M15").Value = -1387.0322
$ws.Range("H33").Value = 162
$ws.Range("I33").Value = 181.38461
$ws.Range("K33").Value = 181.38461
$ws.Range("M33").Value = 47.61538999999999
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1455
$ws.Range("I32").Value = 1091.4043
$ws.Range("K32").Value = 1091.4043
$ws.Range("M32").Value = -804.4042999999999
$ws.Range("H45").Value = 1270
$ws.Range("I45").Value = 783.3333
$ws.Range("K45").Value = 783.3333
$ws.Range("M45").Value = -406.3333
$ws.Range("H61").Value = 3332.8333
$ws.Range("I61").Value = 3000
$ws.Range("J61").Value = 3499.25
$ws.Range("K61").Value = 3000
$ws.Range("L61").Value = 3499.25
$ws.Range("M61").Value = -2788
$ws.Range("N61").Value = -3923.25
$ws.Range("H132").Value = 2626.3333
$ws.Range("I132").Value = 2642.125
$ws.Range("K132").Value = 7926.375
$ws.Range("M132").Value = -5396.375
$ws.Range("H136").Value = 3332.8333
$ws.Range("I136").Value = 3000
$ws.Range("J136").Value = 3499.25
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 10497.75
$ws.Range("M136").Value = -6450
$ws.Range("N136").Value = -15597.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 10000
$ws.Range("J18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("N18").Value = -11058
$ws.Range("H20").Value = 5989.222
$ws.Range("I20").Value = 5237.875
$ws.Range("K20").Value = 5237.875
$ws.Range("M20").Value = -4990.875
$ws.Range("H105").Value = 2000
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -5494
$ws.Range("H107").Value = 805.0833
$ws.Range("I107").Value = 751.5714
$ws.Range("J107").Value = 880
$ws.Range("K107").Value = 751.5714
$ws.Range("L107").Value = 880
$ws.Range("M107").Value = 1168.4286
$ws.Range("N107").Value = -4720
$ws.Range("H134").Value = 15337.333
$ws.Range("I134").Value = 15337.333
$ws.Range("K134").Value = 46011.999
$ws.Range("M134").Value = -43476.999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 615.4
$ws.Range("I16").Value = 556.75
$ws.Range("K16").Value = 556.75
$ws.Range("M16").Value = -269.75
$ws.Range("H31").Value = 2670.0952
$ws.Range("I31").Value = 2237.4443
$ws.Range("K31").Value = 2237.4443
$ws.Range("M31").Value = -1942.4443
$ws.Range("H34").Value = 2670.0952
$ws.Range("I34").Value = 2237.4443
$ws.Range("K34").Value = 2237.4443
$ws.Range("M34").Value = -2035.4443
$ws.Range("H86").Value = 7889.8335
$ws.Range("I86").Value = 6496.3335
$ws.Range("J86").Value = 9283.333000000001
$ws.Range("K86").Value = 6496.3335
$ws.Range("L86").Value = 9283.333000000001
$ws.Range("M86").Value = -5373.3335
$ws.Range("N86").Value = -11529.333
$ws.Range("H89").Value = 7889.8335
$ws.Range("I89").Value = 6496.3335
$ws.Range("J89").Value = 9283.333000000001
$ws.Range("K89").Value = 32481.6675
$ws.Range("L89").Value = 46416.665
$ws.Range("M89").Value = -26865.6675
$ws.Range("N89").Value = -57648.665
$ws.Range("H105").Value = 508.2857
$ws.Range("I105").Value = 468
$ws.Range("J105").Value = 750
$ws.Range("K105").Value = 468
$ws.Range("L105").Value = 750
$ws.Range("M105").Value = 1279
$ws.Range("N105").Value = -4244
$ws.Range("H107").Value = 522.73914
$ws.Range("I107").Value = 504.44446
$ws.Range("J107").Value = 588.6
$ws.Range("K107").Value = 504.44446
$ws.Range("L107").Value = 588.6
$ws.Range("M107").Value = 1415.55554
$ws.Range("N107").Value = -4428.6
$ws.Range("H113").Value = 615.4
$ws.Range("I113").Value = 556.75
$ws.Range("K113").Value = 556.75
$ws.Range("M113").Value = 1613.25
$ws.Range("H120").Value = 14000
$ws.Range("J120").Value = 14000
$ws.Range("L120").Value = 14000
$ws.Range("N120").Value = -21258
$ws.Range("H122").Value = 1517.5
$ws.Range("I122").Value = 1330.3334
$ws.Range("J122").Value = 2079
$ws.Range("K122").Value = 3991.0002
$ws.Range("L122").Value = 6237
$ws.Range("M122").Value = -1541.0002
$ws.Range("N122").Value = -11137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 289.2857
$ws.Range("H38").Value = 170.5
$ws.Range("I38").Value = 241
$ws.Range("K38").Value = 723
$ws.Range("M38").Value = -376
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H68").Value = 1032.6666
$ws.Range("I68").Value = 999.3333
$ws.Range("J68").Value = 1066
$ws.Range("K68").Value = 2997.9999
$ws.Range("L68").Value = 3198
$ws.Range("M68").Value = -2186.9999
$ws.Range("N68").Value = -4820
$ws.Range("H71").Value = 1032.6666
$ws.Range("I71").Value = 999.3333
$ws.Range("J71").Value = 1066
$ws.Range("K71").Value = 8993.9997
$ws.Range("L71").Value = 9594
$ws.Range("M71").Value = -4937.9997
$ws.Range("N71").Value = -17706
$ws.Range("H80").Value = 9388.6
$ws.Range("J80").Value = 9611
$ws.Range("L80").Value = 28833
$ws.Range("N80").Value = -30705
$ws.Range("H83").Value = 9388.6
$ws.Range("J83").Value = 9611
$ws.Range("L83").Value = 86499
$ws.Range("N83").Value = -95859
$ws.Range("H107").Value = 1247
$ws.Range("I107").Value = 1495
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 4485
$ws.Range("L107").Value = 2997
$ws.Range("M107").Value = -2565
$ws.Range("N107").Value = -6837
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 30000
$ws.Range("J15").Value = 30000
$ws.Range("L15").Value = 30000
$ws.Range("N15").Value = -30576
$ws.Range("H70").Value = 66676028
$ws.Range("I70").Value = 66676028
$ws.Range("K70").Value = 66676028
$ws.Range("M70").Value = -66675758
$ws.Range("H73").Value = 66676028
$ws.Range("I73").Value = 66676028
$ws.Range("K73").Value = 66676028
$ws.Range("M73").Value = -66675092
$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 30000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -31996
$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -99984
$ws.Range("H97").Value = 1889.3334
$ws.Range("I97").Value = 408.1111
$ws.Range("J97").Value = 6333
$ws.Range("K97").Value = 408.1111
$ws.Range("L97").Value = 6333
$ws.Range("M97").Value = 87.88889999999998
$ws.Range("N97").Value = -7325
$ws.Range("H113").Value = 1337.5834
$ws.Range("I113").Value = 981.625
$ws.Range("K113").Value = 981.625
$ws.Range("M113").Value = 1188.375
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 950.9
$ws.Range("I22").Value = 838.625
$ws.Range("K22").Value = 838.625
$ws.Range("M22").Value = -543.625
$ws.Range("H27").Value = 950.9
$ws.Range("I27").Value = 838.625
$ws.Range("K27").Value = 838.625
$ws.Range("M27").Value = -731.625
$ws.Range("H82").Value = 1584
$ws.Range("I82").Value = 1300.8
$ws.Range("K82").Value = 1300.8
$ws.Range("M82").Value = -939.8
$ws.Range("H85").Value = 1584
$ws.Range("I85").Value = 1300.8
$ws.Range("K85").Value = 1300.8
$ws.Range("M85").Value = -52.79999999999995
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H100").Value = 649.5
$ws.Range("I100").Value = 399
$ws.Range("J100").Value = 900
$ws.Range("K100").Value = 399
$ws.Range("L100").Value = 900
$ws.Range("M100").Value = 142
$ws.Range("N100").Value = -1982
$ws.Range("H132").Value = 2439.8
$ws.Range("I132").Value = 2439.8
$ws.Range("K132").Value = 7319.400000000001
$ws.Range("M132").Value = -4789.400000000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1718.6666
$ws.Range("I107").Value = 1093.6
$ws.Range("K107").Value = 3280.8
$ws.Range("M107").Value = -1360.8
$ws.Range("H132").Value = 3749.5
$ws.Range("I132").Value = 3749.5
$ws.Range("K132").Value = 11248.5
$ws.Range("M132").Value = -8718.5
